$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A45").Value = "2025-04-29 04:52:09"
$ws.Range("B45").Value = 142
